$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.435.38"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "3.527.99"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.45%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("E12").Value = "  +3.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").Value = "4.091.87"
$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").Value = "3.531.00"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("D18").Value = "66.408.42"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("E19").Value = "  +1.88%  "

$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("E22").Value = "  +7.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "640.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.386"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").Value = "3.264.24"
$ws.Range("E42").Value = "  +8.24%  "

$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.14%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0424"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("E49").Value = "  +2.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
